$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Unmerge the old "room_data" header block and clear the old content area ---
$ws.Range("H6:L7").UnMerge()
$ws.Range("C4:N12").ClearContents()
# the two new columns that now sit where the merged header used to be must lose
# the centered style that the old merge left behind
$ws.Range("H6:I7").ClearFormats()

# --- 2. Cells that keep re-using text already present in the workbook
#        (title, unchanged headers/types, and the shifted room_data sub-table) ---
$ws.Range("C4").Value = "Estructura de los datos a almacenar"

$ws.Range("C6").Value = "name"
$ws.Range("G6").Value = "hotel_features"
$ws.Range("H6").Value = "popular_services"
$ws.Range("J6").Value = "room_data"

$ws.Range("C7").Value = "string"
$ws.Range("D7").Value = "string"
$ws.Range("G7").Value = "list"
$ws.Range("H7").Value = "list"
$ws.Range("J7").Value = "dictionary"

$ws.Range("J8").Value = "type_of_room"
$ws.Range("K8").Value = "capacity"
$ws.Range("L8").Value = "price"
$ws.Range("M8").Value = "options"
$ws.Range("N8").Value = "room_features"

$ws.Range("J9").Value = "string"
$ws.Range("K9").Value = "int"
$ws.Range("L9").Value = "int"
$ws.Range("M9").Value = "list"
$ws.Range("N9").Value = "list/dict"

$ws.Range("C10").Value = """Ilunion Bel-Art"""
$ws.Range("G10").Value = "[Vistas a la ciudad, Admite mascotas, …]"
$ws.Range("J10").Value = """Habitación doble"""
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 285
$ws.Range("M10").Value = "[""desayuno"", ""reembolso""]"
$ws.Range("N10").Value = "[""15m2""]"

$ws.Range("N12").Value = "Será un diccionario si puedo extraer el nombre de cada categoria"

# --- 3. Brand-new text, entered in the order the author likely typed it ---
$ws.Range("E6").Value = "hotel_score"
$ws.Range("E7").Value = "double"
# "8.5" is stored as TEXT in the target file, not a number - force text entry
# then drop the text number-format again so no stray style sticks to the cell
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "8.5"
$ws.Range("E10").ClearFormats()

$ws.Range("D6").Value = "address"
$ws.Range("D10").Value = """Calle la que sea"""

$ws.Range("F6").Value = "hotel_description"
$ws.Range("F7").Value = "list of strings"

$ws.Range("I6").Value = "hotel_scores"
$ws.Range("I8").Value = "key for category"
$ws.Range("I7").Value = "dictionary"

# --- 4. Re-create the merged "room_data" header cells with centered style ---
$ws.Range("J7:N7").Merge()
$ws.Range("J6:N6").Merge()

# --- 5. Column widths to roughly match the new layout ---
# (old custom widths belonged to columns whose content has since shifted; re-apply them
#  to the columns that now hold that content, and relax the now-unused column K back
#  toward the sheet's default width)
$ws.Range("K1").ColumnWidth = 9.67
$ws.Range("D1:F1").ColumnWidth = 17.67
$ws.Range("G1").ColumnWidth = 35
$ws.Range("H1:I1").ColumnWidth = 15.33
$ws.Range("J1").ColumnWidth = 16.83
$ws.Range("M1").ColumnWidth = 25.5

# --- 6. Selection as in the saved file ---
$ws.Range("H22").Select()
